# Applies the vic-key-outbreaks data refresh:
#  - Header "Cluster Name" -> "Cluster name"
#  - Cluster list re-sourced/re-sorted alphabetically (renames, additions, removals)
#  - Active cases counts updated; row for "Aspect Autism..." loses its count (blank)
#  - Table grows from A1:B39 to A1:B46 (7 new clusters appended at the bottom of the sort)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Cluster name"

$clusterNames = @(
  "139 Highett St Apartment Complex Richmond",
  "3153 Sacred Heart Community St Kilda Tier 1A",
  "3528 Ottoman Village Aged Care Broadmeadows",
  "3600 Belvedere Aged Care Noble Park",
  "3612 BlueCross Glengowrie",
  "3652 Regis Aged Care Dandenong North",
  "3684 Homestyle Aged Care Langford Grange Cranbourne East",
  "3824 Estia Health South Morang",
  "3980 Arcare Keysborough Aged Care Keysborough",
  "4518 Regis Aged Care Fawkner",
  "ACFS Port Logistics Altona",
  "Alfred Health The Alfred Hospital Melbourne",
  "Armstrong Creek School Armstrong Creek",
  "Aspect Autism Spectrum Australia Disability Service Heatherton",
  "Berwick Fields Primary School Berwick",
  "Berwick Lodge Primary School Berwick",
  "Bridgewood Primary School Officer",
  "Bubup Womindjeka Family and Children's Centre Port Melbourne",
  "Clifton Hill Primary School Clifton Hill",
  "Dandenong North Primary School Dandenong",
  "Elements Childcare Warralily Armstrong Creek",
  "Honeyeater Hairdressers Bendigo",
  "Inghams Chicken Warehouse Lyndhurst",
  "Inghams Enterprises Somerville",
  "JBS Australia Brooklyn",
  "KingKids Early Learning Centre and Kindergarten Hallam",
  "Kmart Distribution Centre Truganina",
  "Lighthouse Christian College Cranbourne",
  "Lilydale Motor Inn Lilydale",
  "Lowanna College Newborough",
  "McQuinns Gym Bendigo",
  "Metcash Limited Distribution Centre Laverton North",
  "Monash Health Casey Hospital Emergency Department Tier 1B",
  "Nido Early School Wyndham Vale",
  "Pelican Childcare Cragieburn",
  "Saint Augustines Primary School Wodonga",
  "St Mary's Primary School Swan Hill",
  "St Vincents Hospital Emergency Department Melbourne",
  "TUROSI PTY LTD Thomastown",
  "The Royal Children's Hospital Parkville",
  "Vizzarri Farms Koo Wee Rup",
  "Werribee Mercy Hospital Emergency Department",
  "Western Health Sunshine Hospital Emergency Department",
  "Wodonga Cemetery Wodonga",
  "Wodonga South Primary School Wodonga"
)

$activeCases = @(
  13,
  13,
  25,
  23,
  33,
  27,
  15,
  52,
  22,
  11,
  14,
  13,
  14,
  $null,
  13,
  19,
  10,
  10,
  15,
  11,
  28,
  12,
  10,
  19,
  14,
  11,
  11,
  11,
  13,
  23,
  24,
  19,
  14,
  15,
  12,
  10,
  14,
  25,
  16,
  11,
  28,
  41,
  21,
  21,
  12
)

for ($i = 0; $i -lt $clusterNames.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $clusterNames[$i]
  if ($activeCases[$i] -eq $null) {
    $ws.Cells.Item($row, 2).ClearContents()
  } else {
    $ws.Cells.Item($row, 2).Value = $activeCases[$i]
  }
}

